$d = $word.ActiveDocument

# Change 1: Split "MétodoApresentação de Relatório do estágio de pesquisa individual."
# into "Método" + line break + "Apresentação de Relatório do estágio de pesquisa individual."
$d.Content.Find.Execute(
    "MétodoApresentação de Relatório do estágio de pesquisa individual.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Método^lApresentação de Relatório do estágio de pesquisa individual.",
    2)

# Change 2: Split the bibliography text into two runs separated by a line break.
$d.Content.Find.Execute(
    "A ser definida pelo supervisor dentro das particularidades da área selecionada.Artigos científicos e técnicos atualizados que tratem dos temas abordados durante o estágio e que forneçam base sólida de conhecimento para o aluno.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A ser definida pelo supervisor dentro das particularidades da área selecionada.^lArtigos científicos e técnicos atualizados que tratem dos temas abordados durante o estágio e que forneçam base sólida de conhecimento para o aluno.",
    2)
